# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ /
# LeveProfitHQ columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
# A handful of rows also gain/lose a lone NQ-or-HQ profit cell (M or N) when
# the recipe's other-quality branch becomes priced/unpriced.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 16797.133
$ws.Range("I40").Value = 4826.3335
$ws.Range("K40").Value = 4826.3335
$ws.Range("M40").Value = -4651.3335
$ws.Range("H55").Value = 160
$ws.Range("I55").Value = 49
$ws.Range("J55").Value = 187.75
$ws.Range("K55").Value = 49
$ws.Range("L55").Value = 187.75
$ws.Range("M55").Value = 165
$ws.Range("N55").Value = -615.75
$ws.Range("H62").Value = 4041.2222
$ws.Range("I62").Value = 3933.875
$ws.Range("K62").Value = 3933.875
$ws.Range("M62").Value = -3309.875
$ws.Range("H64").Value = 8874.875
$ws.Range("I64").Value = 7499.75
$ws.Range("K64").Value = 7499.75
$ws.Range("M64").Value = -7251.75
$ws.Range("H65").Value = 4041.2222
$ws.Range("I65").Value = 3933.875
$ws.Range("K65").Value = 19669.375
$ws.Range("M65").Value = -16549.375
$ws.Range("H67").Value = 8874.875
$ws.Range("I67").Value = 7499.75
$ws.Range("K67").Value = 7499.75
$ws.Range("M67").Value = -6641.75
$ws.Range("H113").Value = 5990.154
$ws.Range("I113").Value = 8123.75
$ws.Range("K113").Value = 8123.75
$ws.Range("M113").Value = -4869.75
$ws.Range("H137").Value = 1589.2325
$ws.Range("I137").Value = 1265.75
$ws.Range("J137").Value = 2193.0667
$ws.Range("K137").Value = 3797.25
$ws.Range("L137").Value = 6579.2001
$ws.Range("M137").Value = -1247.25
$ws.Range("N137").Value = -11679.2001
$ws.Range("H138").Value = 4041
$ws.Range("I138").Value = 2357.6
$ws.Range("J138").Value = 4567.0625
$ws.Range("K138").Value = 7072.799999999999
$ws.Range("L138").Value = 13701.1875
$ws.Range("M138").Value = -1932.799999999999
$ws.Range("N138").Value = -23981.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2544.4
$ws.Range("I5").Value = 2358.2144
$ws.Range("J5").Value = 2978.8333
$ws.Range("K5").Value = 2358.2144
$ws.Range("L5").Value = 2978.8333
$ws.Range("M5").Value = -2246.2144
$ws.Range("N5").Value = -3202.8333
$ws.Range("H98").Value = 36902.668
$ws.Range("J98").Value = 36902.668
$ws.Range("L98").Value = 36902.668
$ws.Range("N98").Value = -42892.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2544.4
$ws.Range("I4").Value = 2358.2144
$ws.Range("J4").Value = 2978.8333
$ws.Range("K4").Value = 2358.2144
$ws.Range("L4").Value = 2978.8333
$ws.Range("M4").Value = -2243.2144
$ws.Range("N4").Value = -3208.8333
$ws.Range("H86").Value = 3912.25
$ws.Range("I86").Value = 2603.2144
$ws.Range("J86").Value = 5744.9
$ws.Range("K86").Value = 2603.2144
$ws.Range("L86").Value = 5744.9
$ws.Range("M86").Value = -1480.2144
$ws.Range("N86").Value = -7990.9
$ws.Range("H89").Value = 3912.25
$ws.Range("I89").Value = 2603.2144
$ws.Range("J89").Value = 5744.9
$ws.Range("K89").Value = 13016.072
$ws.Range("L89").Value = 28724.5
$ws.Range("M89").Value = -7400.072
$ws.Range("N89").Value = -39956.5
$ws.Range("H94").Value = 1187.4
$ws.Range("I94").Value = 1187.4
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1187.4
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -736.4000000000001
$ws.Range("N94").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9692.058999999999
$ws.Range("I16").Value = 718.8889
$ws.Range("K16").Value = 718.8889
$ws.Range("M16").Value = -431.8889
$ws.Range("H31").Value = 5127.0557
$ws.Range("I31").Value = 6000
$ws.Range("J31").Value = 2071.75
$ws.Range("K31").Value = 6000
$ws.Range("L31").Value = 2071.75
$ws.Range("M31").Value = -5705
$ws.Range("N31").Value = -2661.75
$ws.Range("H34").Value = 5127.0557
$ws.Range("I34").Value = 6000
$ws.Range("J34").Value = 2071.75
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 2071.75
$ws.Range("M34").Value = -5798
$ws.Range("N34").Value = -2475.75
$ws.Range("H43").Value = 47642.11
$ws.Range("J43").Value = 49825.855
$ws.Range("L43").Value = 49825.855
$ws.Range("N43").Value = -50193.855
$ws.Range("H96").Value = 39117.4
$ws.Range("J96").Value = 39117.4
$ws.Range("L96").Value = 39117.4
$ws.Range("N96").Value = -44609.4
$ws.Range("H101").Value = 47642.11
$ws.Range("J101").Value = 49825.855
$ws.Range("L101").Value = 49825.855
$ws.Range("N101").Value = -56315.855
$ws.Range("H113").Value = 9692.058999999999
$ws.Range("I113").Value = 718.8889
$ws.Range("K113").Value = 718.8889
$ws.Range("M113").Value = 1451.1111
$ws.Range("H141").Value = 62557
$ws.Range("J141").Value = 64683.3
$ws.Range("L141").Value = 64683.3
$ws.Range("N141").Value = -75043.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2130.6
$ws.Range("I131").Value = 892.4286
$ws.Range("K131").Value = 2677.2858
$ws.Range("M131").Value = 2362.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1859.2727
$ws.Range("I102").Value = 1815.2
$ws.Range("K102").Value = 1815.2
$ws.Range("M102").Value = -193.2
$ws.Range("H113").Value = 6695
$ws.Range("I113").Value = 7901.25
$ws.Range("K113").Value = 7901.25
$ws.Range("M113").Value = -5731.25
$ws.Range("H122").Value = 4700.7617
$ws.Range("I122").Value = 3792.4375
$ws.Range("J122").Value = 7607.4
$ws.Range("K122").Value = 11377.3125
$ws.Range("L122").Value = 22822.2
$ws.Range("M122").Value = -8927.3125
$ws.Range("N122").Value = -27722.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3738.0625
$ws.Range("J46").Value = 3671.818
$ws.Range("L46").Value = 3671.818
$ws.Range("N46").Value = -4047.818
$ws.Range("H55").Value = 347.46667
$ws.Range("I55").Value = 348.16666
$ws.Range("J55").Value = 344.66666
$ws.Range("K55").Value = 348.16666
$ws.Range("L55").Value = 344.66666
$ws.Range("M55").Value = -175.16666
$ws.Range("N55").Value = -690.66666
$ws.Range("H61").Value = 5119.6
$ws.Range("I61").Value = 5149.5
$ws.Range("K61").Value = 5149.5
$ws.Range("M61").Value = -4947.5
$ws.Range("H64").Value = 14000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 14000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 14000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -14450
$ws.Range("H67").Value = 14000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 14000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 14000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -15560
$ws.Range("H113").Value = 5119.6
$ws.Range("I113").Value = 5149.5
$ws.Range("K113").Value = 5149.5
$ws.Range("M113").Value = -2979.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 40000
$ws.Range("I63").Value = 40000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 40000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -39376
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 40000
$ws.Range("I66").Value = 40000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 120000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -116880
$ws.Range("N66").ClearContents()
$ws.Range("H81").Value = 3659.75
$ws.Range("I81").Value = 2754
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 5508
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -4447
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 3659.75
$ws.Range("I84").Value = 2754
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 27540
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -22236
$ws.Range("N84").Value = -110608
$ws.Range("H97").Value = 46997.668
$ws.Range("J97").Value = 46997.668
$ws.Range("L97").Value = 46997.668
$ws.Range("N97").Value = -48979.668
$ws.Range("H122").Value = 4543
$ws.Range("I122").Value = 2957
$ws.Range("K122").Value = 8871
$ws.Range("M122").Value = -6421
